$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new rows (49-60) to the table: cloud service models (IaaS/PaaS/SaaS),
# deployment models (Private/Public/Hybrid/Community) and VDE/VDI.
# Values are written in the same cell order the original author used so the
# shared-string table indices line up with the source workbook.

$ws.Range("A49").Value = "Infrastructure as a Service (IAAS)"
$ws.Range("C49").Value = "A remote company will rent you an entire infrastructure (such as Amazon Web Services)"
$ws.Range("B49").Value = "Definition"

$ws.Range("B50").Value = "Capabilities"
$ws.Range("C50").Value = "Can use remote desktop, filter access by IP, setup virtual machine. Even a website!"

$ws.Range("A51").Value = "Platform as a service (PAAS)"
$ws.Range("B51").Value = "Defintion"

$ws.Range("B52").Value = "Addons"
$ws.Range("C52").Value = "Capable to load addons with a single click"

$ws.Range("C51").Value = "Similar to IAAS, but detatched from physical hardware. (e.g. Heroku). Everything is obfuscated"

$ws.Range("B53").Value = "Advantages"
$ws.Range("C53").Value = "Enables you to get software running live on the internet very quickly"

$ws.Range("A54").Value = "Software as a service (SAAS)"
$ws.Range("C54").Value = "Gets rid of physical media. E.g. Microsoft Office. Software subscriptions, e.g. Dropbox and Google Docs"
$ws.Range("B54").Value = "Definition"

$ws.Range("A55").Value = "Deployment Models"
$ws.Range("B55").Value = "Private Cloud"
$ws.Range("C55").Value = "Soley for use within organisation"

$ws.Range("B56").Value = "Public Cloud"
$ws.Range("C56").Value = "E.g. Microsoft Azure. Anyone can use these."

$ws.Range("B57").Value = "Hybrid Cloud"
$ws.Range("C57").Value = "Some of the cloud is private, some of the cloud is public"

$ws.Range("B58").Value = "Community Cloud"
$ws.Range("C58").Value = '"Members Only" used within groups of users to save capital'

$ws.Range("B59").Value = "Virtual Desktop Environment (VDE)"
$ws.Range("B60").Value = "Virtual Desktop Interface (VDI)"
$ws.Range("C60").Value = "The actual virtualized environment on the cloud"
$ws.Range("C59").Value = "Accessing a remote physical desktop"

$ws.Range("C59").Select()
